$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "301KG0KX9CLV8GLA6QPGKOCZD972HG"
$ws.Range("C2").Value = "At the evening offering I arose up from my humiliation, even with my garment and my robe torn; and I fell on my knees, and spread out my hands to Yahweh my God;"
$ws.Range("D2").Value = "garment"
$ws.Range("E2").Value = "easy"
$ws.Range("F2").Value = "0.01 - 0.25"
$ws.Range("G2").Value = 0.2499896965645608
$ws.Range("H2").Value = 0.21875
$ws.Range("I2").Value = "easy"
$ws.Range("J2").Value = "Si"
$ws.Range("K2").Value = 0.1233
$ws.Range("L2").Value = 0.0261
$ws.Range("M2").Value = 0.1615
$ws.Range("N2").Value = 0.1274
$ws.Range("O2").Value = 0.4622
$ws.Range("P2").Value = 0.3333
$ws.Range("Q2").Value = "neutral:0.74%"
$ws.Range("R2").Value = "Easy:0.01%"
$ws.Range("S2").Value = "easy:0.0%"
$ws.Range("T2").Value = ":0.0%"
$ws.Range("U2").Value = "None"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "301KG0KX9CLV8GLA6QPGKOCZD9X2H6"
$ws.Range("C3").Value = "It shall be to them as a false divination in their sight, who have sworn oaths to them; but he brings iniquity to memory, that they may be taken."
$ws.Range("D3").Value = "divination"
$ws.Range("E3").Value = "neutral"
$ws.Range("F3").Value = "0.26 - 0.5"
$ws.Range("G3").Value = 0.375
$ws.Range("H3").Value = 0.638888918
$ws.Range("I3").Value = "difficult"
$ws.Range("J3").Value = "No"
$ws.Range("K3").Value = 0.1233
$ws.Range("L3").Value = 0.0261
$ws.Range("M3").Value = 0.1615
$ws.Range("N3").Value = 0.1274
$ws.Range("O3").Value = 0.4622
$ws.Range("P3").Value = 0.3333
$ws.Range("Q3").Value = "neutral:84.37%"
$ws.Range("R3").Value = "difficult:10.64%"
$ws.Range("S3").Value = "easy:4.89%"
$ws.Range("T3").Value = "Neutral:0.07%"
$ws.Range("U3").Value = "None"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "301KG0KX9CLV8GLA6QPGKOCZDB2H2U"
$ws.Range("C4").Value = "You have received gifts among men, yes, among the rebellious also, that Yah God might dwell there."
$ws.Range("D4").Value = "gifts"
$ws.Range("E4").Value = "easy"
$ws.Range("F4").Value = "0.01 - 0.25"
$ws.Range("G4").Value = 0.2499889281109571
$ws.Range("H4").Value = 0.089285682
$ws.Range("I4").Value = "easy"
$ws.Range("J4").Value = "Si"
$ws.Range("K4").Value = 0.1233
$ws.Range("L4").Value = 0.0261
$ws.Range("M4").Value = 0.1615
$ws.Range("N4").Value = 0.1274
$ws.Range("O4").Value = 0.4622
$ws.Range("P4").Value = 0.3333
$ws.Range("Q4").Value = "neutral:0.93%"
$ws.Range("R4").Value = "Easy:0.01%"
$ws.Range("S4").Value = "very:0.01%"
$ws.Range("T4").Value = "easy:0.0%"
$ws.Range("U4").Value = "None"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "301KG0KX9CLV8GLA6QPGKOCZDBWH2O"
$ws.Range("C5").Value = "Therefore he poured the fierceness of his anger on him, and the strength of battle; and it set him on fire all around, but he didn't know; and it burned him, but he didn't take it to heart.`""
$ws.Range("D5").Value = "strength"
$ws.Range("E5").Value = "easy"
$ws.Range("F5").Value = "0.01 - 0.25"
$ws.Range("G5").Value = 0.2499878463276225
$ws.Range("H5").Value = 0.166666648
$ws.Range("I5").Value = "easy"
$ws.Range("J5").Value = "Si"
$ws.Range("K5").Value = 0.1233
$ws.Range("L5").Value = 0.0261
$ws.Range("M5").Value = 0.1615
$ws.Range("N5").Value = 0.1274
$ws.Range("O5").Value = 0.4622
$ws.Range("P5").Value = 0.3333
$ws.Range("Q5").Value = "neutral:14.04%"
$ws.Range("R5").Value = "Easy:0.02%"
$ws.Range("S5").Value = "Neutral:0.01%"
$ws.Range("T5").Value = "easy:0.0%"
$ws.Range("U5").Value = "None"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "301KG0KX9CLV8GLA6QPGKOCZDBX2HA"
$ws.Range("C6").Value = "The seventh angel sounded, and great voices in heaven followed, saying, `"The kingdom of the world has become the Kingdom of our Lord, and of his Christ."
$ws.Range("D6").Value = "voices"
$ws.Range("E6").Value = "easy"
$ws.Range("F6").Value = "0.01 - 0.25"
$ws.Range("G6").Value = 0.2499884037545524
$ws.Range("H6").Value = 0.18421044
$ws.Range("I6").Value = "easy"
$ws.Range("J6").Value = "Si"
$ws.Range("K6").Value = 0.1233
$ws.Range("L6").Value = 0.0261
$ws.Range("M6").Value = 0.1615
$ws.Range("N6").Value = 0.1274
$ws.Range("O6").Value = 0.4622
$ws.Range("P6").Value = 0.3333
$ws.Range("Q6").Value = "neutral:1.48%"
$ws.Range("R6").Value = "Easy:0.02%"
$ws.Range("S6").Value = "easy:0.0%"
$ws.Range("T6").Value = "Neutral:0.0%"
$ws.Range("U6").Value = "None"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "302OLP89DZ7TWB5YXD4UFFYHC58AC0"
$ws.Range("C7").Value = "who by the mouth of your servant, David, said, 'Why do the nations rage, and the peoples plot a vain thing?"
$ws.Range("D7").Value = "rage"
$ws.Range("E7").Value = "easy"
$ws.Range("F7").Value = "0.01 - 0.25"
$ws.Range("G7").Value = 0.2499853649494029
$ws.Range("H7").Value = 0.214285633
$ws.Range("I7").Value = "easy"
$ws.Range("J7").Value = "Si"
$ws.Range("K7").Value = 0.1233
$ws.Range("L7").Value = 0.0261
$ws.Range("M7").Value = 0.1615
$ws.Range("N7").Value = 0.1274
$ws.Range("O7").Value = 0.4622
$ws.Range("P7").Value = 0.3333
$ws.Range("Q7").Value = "neutral:2.84%"
$ws.Range("R7").Value = "Easy:0.02%"
$ws.Range("S7").Value = "Neutral:0.01%"
$ws.Range("T7").Value = "easy:0.01%"
$ws.Range("U7").Value = "None"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "302OLP89DZ7TWB5YXD4UFFYHC5KCAE"
$ws.Range("C8").Value = "''You shall keep my Sabbaths, and reverence my sanctuary; I am Yahweh."
$ws.Range("D8").Value = "Sabbaths"
$ws.Range("E8").Value = "easy"
$ws.Range("F8").Value = "0.01 - 0.25"
$ws.Range("G8").Value = 0.2499822299869348
$ws.Range("H8").Value = 0.34375
$ws.Range("I8").Value = "neutral"
$ws.Range("J8").Value = "No"
$ws.Range("K8").Value = 0.1233
$ws.Range("L8").Value = 0.0261
$ws.Range("M8").Value = 0.1615
$ws.Range("N8").Value = 0.1274
$ws.Range("O8").Value = 0.4622
$ws.Range("P8").Value = 0.3333
$ws.Range("Q8").Value = "neutral:1.09%"
$ws.Range("R8").Value = "Easy:0.03%"
$ws.Range("S8").Value = "easy:0.01%"
$ws.Range("T8").Value = "Neutral:0.01%"
$ws.Range("U8").Value = "None"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "302OLP89DZ7TWB5YXD4UFFYHC7KACG"
$ws.Range("C9").Value = "Neither have we inheritance in the son of Jesse!"
$ws.Range("D9").Value = "Jesse"
$ws.Range("E9").Value = "easy"
$ws.Range("F9").Value = "0.01 - 0.25"
$ws.Range("G9").Value = 0.009698101353405764
$ws.Range("H9").Value = 0.323529285
$ws.Range("I9").Value = "neutral"
$ws.Range("J9").Value = "No"
$ws.Range("K9").Value = 0.1233
$ws.Range("L9").Value = 0.0261
$ws.Range("M9").Value = 0.1615
$ws.Range("N9").Value = 0.1274
$ws.Range("O9").Value = 0.4622
$ws.Range("P9").Value = 0.3333
$ws.Range("Q9").Value = "easy:96.12%"
$ws.Range("R9").Value = "neutral:3.76%"
$ws.Range("S9").Value = "very:0.08%"
$ws.Range("T9").Value = "Easy:0.02%"
$ws.Range("U9").Value = "Neutral:0.01%"

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "302OLP89DZ7TWB5YXD4UFFYHC7PCAN"
$ws.Range("C10").Value = "Only the firstborn among animals, which is made a firstborn to Yahweh, no man may dedicate it; whether an ox or sheep, it is Yahweh's."
$ws.Range("D10").Value = "animals"
$ws.Range("E10").Value = "easy"
$ws.Range("F10").Value = "0.01 - 0.25"
$ws.Range("G10").Value = 0.2499928275566013
$ws.Range("H10").Value = 0.222222222
$ws.Range("I10").Value = "easy"
$ws.Range("J10").Value = "Si"
$ws.Range("K10").Value = 0.1233
$ws.Range("L10").Value = 0.0261
$ws.Range("M10").Value = 0.1615
$ws.Range("N10").Value = 0.1274
$ws.Range("O10").Value = 0.4622
$ws.Range("P10").Value = 0.3333
$ws.Range("Q10").Value = "neutral:0.49%"
$ws.Range("R10").Value = "very:0.03%"
$ws.Range("S10").Value = "Easy:0.01%"
$ws.Range("T10").Value = "easy:0.0%"
$ws.Range("U10").Value = "None"

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "302U8RURJZ1WF35NXY44RD66WJ8NVH"
$ws.Range("C11").Value = "Don't be desirous of his dainties, since they are deceitful food."
$ws.Range("D11").Value = "dainties"
$ws.Range("E11").Value = "neutral"
$ws.Range("F11").Value = "0.26 - 0.5"
$ws.Range("G11").Value = 0.375
$ws.Range("H11").Value = 0.638888918
$ws.Range("I11").Value = "difficult"
$ws.Range("J11").Value = "No"
$ws.Range("K11").Value = 0.1233
$ws.Range("L11").Value = 0.0261
$ws.Range("M11").Value = 0.1615
$ws.Range("N11").Value = 0.1274
$ws.Range("O11").Value = 0.4622
$ws.Range("P11").Value = 0.3333
$ws.Range("Q11").Value = "neutral:81.26%"
$ws.Range("R11").Value = "easy:17.86%"
$ws.Range("S11").Value = "difficult:0.81%"
$ws.Range("T11").Value = "Neutral:0.06%"
$ws.Range("U11").Value = "None"

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "302U8RURJZ1WF35NXY44RD66WJ8VNP"
$ws.Range("C12").Value = "The floods have lifted up, Yahweh, the floods have lifted up their voice."
$ws.Range("D12").Value = "voice"
$ws.Range("E12").Value = "easy"
$ws.Range("F12").Value = "0.01 - 0.25"
$ws.Range("G12").Value = 0.249989758317862
$ws.Range("H12").Value = 0.266666692
$ws.Range("I12").Value = "neutral"
$ws.Range("J12").Value = "No"
$ws.Range("K12").Value = 0.1233
$ws.Range("L12").Value = 0.0261
$ws.Range("M12").Value = 0.1615
$ws.Range("N12").Value = 0.1274
$ws.Range("O12").Value = 0.4622
$ws.Range("P12").Value = 0.3333
$ws.Range("Q12").Value = "neutral:0.75%"
$ws.Range("R12").Value = "very:0.22%"
$ws.Range("S12").Value = "Easy:0.02%"
$ws.Range("T12").Value = "easy:0.0%"
$ws.Range("U12").Value = "None"
